$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.282.23"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "3.020.35"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "518.29"
$ws.Range("E5").Value = "  +4.85%  "

$ws.Range("D6").Value = "141.66"
$ws.Range("E6").Value = "  +5.32%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "0.438"
$ws.Range("E8").Value = "  +3.50%  "

$ws.Range("D9").Value = "7.58"
$ws.Range("E9").Value = "  +5.15%  "

$ws.Range("E10").Value = "  +5.75%  "

$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("E12").Value = "  +2.20%  "

$ws.Range("D13").Value = "3.543.93"
$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("D14").Value = "26.07"
$ws.Range("E14").Value = "  +4.56%  "

$ws.Range("D15").Value = "0.0000160"
$ws.Range("E15").Value = "  +10.94%  "

$ws.Range("D16").Value = "57.308.13"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("D17").Value = "3.025.40"
$ws.Range("E17").Value = "  +1.27%  "

$ws.Range("D18").Value = "6.00"
$ws.Range("E18").Value = "  +2.50%  "

$ws.Range("D19").Value = "12.69"
$ws.Range("E19").Value = "  +3.14%  "

$ws.Range("D20").Value = "7.95"
$ws.Range("E20").Value = "  +3.21%  "

$ws.Range("D21").Value = "331.05"
$ws.Range("E21").Value = "  +2.63%  "

$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").Value = "0.488"
$ws.Range("E23").Value = "  +5.34%  "

$ws.Range("D24").Value = "64.25"
$ws.Range("E24").Value = "  +5.11%  "

$ws.Range("E25").Value = "  +6.14%  "

$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("D27").Value = "0.0₃0924"
$ws.Range("E27").Value = "  +5.42%  "

$ws.Range("D28").Value = "6.78"
$ws.Range("E28").Value = "  +3.21%  "

$ws.Range("D29").Value = "7.21"
$ws.Range("E29").Value = "  +6.75%  "

$ws.Range("E30").Value = "  +6.42%  "

$ws.Range("E31").Value = "  +3.98%  "

$ws.Range("D32").Value = "20.78"
$ws.Range("E32").Value = "  +4.54%  "

$ws.Range("D33").Value = "158.69"
$ws.Range("E33").Value = "  +5.40%  "

$ws.Range("D34").Value = "4.65"
$ws.Range("E34").Value = "  +3.63%  "

$ws.Range("D35").Value = "5.77"
$ws.Range("E35").Value = "  +2.46%  "

$ws.Range("D36").Value = "1.30"
$ws.Range("E36").Value = "  +1.66%  "

$ws.Range("D37").Value = "24.45"
$ws.Range("E37").Value = "  +4.29%  "

$ws.Range("D38").Value = "0.0678"
$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("D39").Value = "3.055.91"
$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("D40").Value = "37.45"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").Value = "0.656"
$ws.Range("E42").Value = "  +2.63%  "

$ws.Range("D43").Value = "2.303.76"
$ws.Range("E43").Value = "  +6.01%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.45"
$ws.Range("E44").Value = "  +1.97%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "3.74"
$ws.Range("E45").Value = "  +5.26%  "

$ws.Range("D46").Value = "1.01"
$ws.Range("E46").Value = "  -0.19%  "

$ws.Range("D47").Value = "2.02"
$ws.Range("E47").Value = "  +8.12%  "

$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("D49").Value = "5.91"
$ws.Range("E49").Value = "  +5.86%  "

$ws.Range("D50").Value = "19.51"
$ws.Range("E50").Value = "  +1.09%  "

$ws.Range("D51").Value = "0.0882"
$ws.Range("E51").Value = "  +3.85%  "
